# Updated cryptos list on Mon Feb 12 09:50:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.948.83'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.481.96'
$ws.Range("E3").Value = '  -1.60%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '316.36'
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("D6").Value = '105.14'
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("E7").Value = '  -3.01%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  -3.61%  '
$ws.Range("D10").Value = '38.78'
$ws.Range("E10").Value = '  -5.19%  '
$ws.Range("D11").Value = '20.11'
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("E12").Value = '  -3.20%  '
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").Value = '2.870.46'
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("D16").Value = '2.479.80'
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("E17").Value = '  -3.66%  '
$ws.Range("D18").Value = '47.857.17'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("E19").Value = '  +8.76%  '
$ws.Range("E20").Value = '  -4.20%  '
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("D22").Value = '0.0₃0928'
$ws.Range("E22").Value = '  -2.15%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '272.74'
$ws.Range("E23").Value = '  +2.76%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '70.70'
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("E28").Value = '  +2.43%  '
$ws.Range("E29").Value = '  -4.94%  '
$ws.Range("E30").Value = '  -4.38%  '
$ws.Range("D31").Value = '34.42'
$ws.Range("E31").Value = '  -4.44%  '
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("E34").Value = '  -5.32%  '
$ws.Range("E35").Value = '  -2.64%  '
$ws.Range("E36").Value = '  -2.88%  '
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("D38").Value = '4.52'
$ws.Range("E38").Value = '  -4.52%  '
$ws.Range("E39").Value = '  -4.95%  '
$ws.Range("D40").Value = '122.50'
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("E42").Value = '  +0.61%  '
$ws.Range("D43").Value = '21.98'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '1.996.82'
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("E48").Value = '  -2.37%  '
$ws.Range("D49").Value = '8.88'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("D51").Value = '78.30'
$ws.Range("E51").Value = '  -1.43%  '
